$d = $word.ActiveDocument

# "Bradley Kersting" row, last column (R6) currently holds a single
# empty paragraph/run. Fill it with "sys_call" and add a second
# paragraph "I/O Scheduler" below it, inheriting the cell's existing
# run/paragraph formatting (arial, sz 20).
$table = $d.Tables.Item(1)
$cell = $table.Cell(4, 6)
$cell.Range.InsertAfter("sys_call`rI/O Scheduler")
